$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.868.61"
$ws.Range("E2").Value = "'  -1.44%  "
$ws.Range("D3").Value = "'2.357.28"
$ws.Range("E3").Value = "'  -2.36%  "
$ws.Range("E4").Value = "'  -0.20%  "
$ws.Range("D5").Value = "'322.08"
$ws.Range("E5").Value = "'  -0.42%  "
$ws.Range("D6").Value = "'106.25"
$ws.Range("E6").Value = "'  +1.47%  "
$ws.Range("D7").Value = "'0.639"
$ws.Range("E7").Value = "'  -1.03%  "
$ws.Range("D9").Value = "'0.619"
$ws.Range("E9").Value = "'  -6.06%  "
$ws.Range("D10").Value = "'41.14"
$ws.Range("E10").Value = "'  -2.37%  "
$ws.Range("D11").Value = "'0.0926"
$ws.Range("E11").Value = "'  -2.59%  "
$ws.Range("D12").Value = "'8.48"
$ws.Range("E12").Value = "'  -1.81%  "
$ws.Range("D13").Value = "'0.998"
$ws.Range("E13").Value = "'  -4.49%  "
$ws.Range("E14").Value = "'  -0.06%  "
$ws.Range("D15").Value = "'16.07"
$ws.Range("E15").Value = "'  -7.26%  "
$ws.Range("D16").Value = "'2.712.45"
$ws.Range("E16").Value = "'  -2.40%  "
$ws.Range("D17").Value = "'2.338.01"
$ws.Range("E17").Value = "'  -2.98%  "
$ws.Range("D18").Value = "'42.837.82"
$ws.Range("E18").Value = "'  -1.59%  "
$ws.Range("D19").Value = "'7.76"
$ws.Range("E19").Value = "'  +3.98%  "
$ws.Range("E20").Value = "'  -3.30%  "
$ws.Range("D21").Value = "'77.35"
$ws.Range("E21").Value = "'  +2.44%  "
$ws.Range("D22").Value = "'3.60"
$ws.Range("E22").Value = "'  +3.13%  "
$ws.Range("D23").Value = "'261.17"
$ws.Range("E23").Value = "'  -0.06%  "
$ws.Range("D24").Value = "'2.34"
$ws.Range("E24").Value = "'  -4.37%  "
$ws.Range("D25").Value = "'9.66"
$ws.Range("E25").Value = "'  -0.48%  "
$ws.Range("E26").Value = "'  +0.17%  "
$ws.Range("D27").Value = "'11.46"
$ws.Range("E27").Value = "'  -3.81%  "
$ws.Range("D28").Value = "'23.36"
$ws.Range("E28").Value = "'  +2.35%  "
$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = "'  -0.66%  "
$ws.Range("D30").Value = "'175.01"
$ws.Range("E30").Value = "'  -2.49%  "
$ws.Range("D31").Value = "'36.50"
$ws.Range("E31").Value = "'  -4.12%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.0894"
$ws.Range("E32").Value = "'  -4.51%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'3.00"
$ws.Range("E33").Value = "'  -7.17%  "
$ws.Range("D34").Value = "'6.11"
$ws.Range("E34").Value = "'  +2.55%  "
$ws.Range("E35").Value = "'  -1.47%  "
$ws.Range("E36").Value = "'  +6.26%  "
$ws.Range("D37").Value = "'4.64"
$ws.Range("E37").Value = "'  -5.38%  "
$ws.Range("D38").Value = "'0.0360"
$ws.Range("E38").Value = "'  -3.14%  "
$ws.Range("D39").Value = "'3.82"
$ws.Range("E39").Value = "'  -4.00%  "
$ws.Range("D40").Value = "'2.69"
$ws.Range("E40").Value = "'  -7.00%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "'71.74"
$ws.Range("E41").Value = "'  +2.48%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.47"
$ws.Range("E42").Value = "'  -9.83%  "
$ws.Range("D43").Value = "'0.233"
$ws.Range("E43").Value = "'  -0.17%  "
$ws.Range("E44").Value = "'  -0.19%  "
$ws.Range("D45").Value = "'115.23"
$ws.Range("E45").Value = "'  -9.32%  "
$ws.Range("D46").Value = "'11.91"
$ws.Range("E46").Value = "'  -5.69%  "
$ws.Range("D47").Value = "'5.53"
$ws.Range("E47").Value = "'  -2.70%  "
$ws.Range("D48").Value = "'9.24"
$ws.Range("E48").Value = "'  -4.02%  "
$ws.Range("D49").Value = "'85.02"
$ws.Range("E49").Value = "'  +7.01%  "
$ws.Range("D50").Value = "'73.89"
$ws.Range("E50").Value = "'  +0.33%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.26"
$ws.Range("E51").Value = "'  -4.48%  "
